$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "IVA, cashback, 'imposto do pecado': a reforma em 7 pontos"
$ws.Range("D2").Value = "https://g1.globo.com/politica/noticia/2023/12/15/entenda-em-7-pontos-a-reforma-tributaria.ghtml"
$ws.Range("E2").Value = "20/12/2023 20:17:33"

# Row 3
$ws.Range("C3").Value = "Senado aprova medida que muda regras de descontos do ICMS "
$ws.Range("D3").Value = "https://g1.globo.com/politica/noticia/2023/12/20/senado-aprova-texto-base-de-mp-que-muda-regras-de-descontos-do-icms-para-aumentar-arrecadacao-federal.ghtml"
$ws.Range("E3").Value = "20/12/2023 20:17:33"

# Row 4
$ws.Range("C4").Value = "STJ quebra sigilos de Cláudio Castro após operação contra irmão"
$ws.Range("D4").Value = "https://g1.globo.com/rj/rio-de-janeiro/noticia/2023/12/20/stj-quebra-sigilos-bancarios-e-telefonico-de-claudio-castro-e-pf-encontra-na-casa-do-irmao-dinheiro-em-caixa-de-remedios.ghtml"
$ws.Range("E4").Value = "20/12/2023 20:17:33"

# Row 5
$ws.Range("C5").Value = "Deputado do PT bate na cara de colega e usa termo homofóbico contra Nikolas"
$ws.Range("D5").Value = "https://noticias.uol.com.br/politica/ultimas-noticias/2023/12/20/vice-presidente-do-pt-da-tapa-na-cara-de-deputado-durante-sessao.htm"
$ws.Range("E5").Value = "20/12/2023 20:17:33"

# Row 6
$ws.Range("C6").Value = "Senado aprova por 48 votos a 22 proposta de Haddad para elevar arrecadação"
$ws.Range("D6").Value = "https://economia.uol.com.br/noticias/redacao/2023/12/20/senado-mp-subvencao-beneficios-fiscais.htm"
$ws.Range("E6").Value = "20/12/2023 20:17:33"

# Row 7
$ws.Range("C7").Value = "Apesar de ameaças de Milei, argentinos protestam contra novo governo"
$ws.Range("D7").Value = "https://noticias.uol.com.br/internacional/ultimas-noticias/2023/12/20/protesto-contra-milei-argentina-nas-ruas-apos-ameacas.htm"
$ws.Range("E7").Value = "20/12/2023 20:17:33"

# Row 8 (only date changes)
$ws.Range("E8").Value = "20/12/2023 20:17:33"

# Row 9 (only date changes)
$ws.Range("E9").Value = "20/12/2023 20:17:33"

# Row 10 (only date changes)
$ws.Range("E10").Value = "20/12/2023 20:17:33"

# Row 11
$ws.Range("C11").Value = "Regulamentação de cigarros eletrônicos está em discussão no Brasil"
$ws.Range("D11").Value = "https://www.cnnbrasil.com.br/nacional/regulamentacao-de-cigarros-eletronicos-esta-em-discussao-no-brasil/"
$ws.Range("E11").Value = "20/12/2023 20:17:33"

# Row 12
$ws.Range("C12").Value = "Lula deve passar Réveillon no litoral do Rio de Janeiro"
$ws.Range("D12").Value = "https://www.cnnbrasil.com.br/politica/lula-deve-passar-reveillon-no-litoral-do-rio-de-janeiro/"
$ws.Range("E12").Value = "20/12/2023 20:17:33"

# Row 13
$ws.Range("C13").Value = "STJ autoriza quebra de sigilos fiscal e telemático de Cláudio Castro"
$ws.Range("D13").Value = "https://www.cnnbrasil.com.br/politica/stj-autoriza-quebra-de-sigilos-fiscal-e-telematico-de-claudio-castro/"
$ws.Range("E13").Value = "20/12/2023 20:17:33"
